$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Replicate row 9's formatting onto the new rows 10-12 so the new cells ---
# --- pick up the same cell styles (s="1"/"0"/"2") used by the existing data rows.
$ws.Range("A9:H9").Copy()
$ws.Range("A10:H10").PasteSpecial(-4122)

$ws.Range("A9:F9").Copy()
$ws.Range("A11:F11").PasteSpecial(-4122)

$ws.Range("C9:D9").Copy()
$ws.Range("C12:D12").PasteSpecial(-4122)

# --- Row 10: new review for "com.hamxa.shaynachim" / "bitcoin free" ---
$ws.Range("A10").Value = "com.hamxa.shaynachim"
$ws.Range("B10").Value = "bitcoin free"
$ws.Range("C10").Value = "hermanliran@gmail.com"
$ws.Range("D10").Value = "nevilgreen@gmail.com"
$ws.Range("E10").Value = "27/5/2019 15:59"
$ws.Range("F10").Value = "bitcoin app especially for beginners with all the info needed. Just great. Keep on the great guide."
$ws.Range("G10").Value = "no"

# --- Row 11: second new review for the same appid/keyword ---
$ws.Range("A11").Value = "com.hamxa.shaynachim"
$ws.Range("B11").Value = "bitcoin free"
$ws.Range("C11").Value = "shmualtamara@gmail.com"
$ws.Range("D11").Value = "shmulmaor2@gmail.com"
$ws.Range("E11").Value = "27/5/2019 15:59"
$ws.Range("F11").Value = "it is exactly what I searched for such a long time. This app saved me a lot of time"

# --- Hyperlinks for the new e-mail addresses ---
$ws.Hyperlinks.Add($ws.Range("C10"), "mailto:hermanliran@gmail.com", [Type]::Missing, [Type]::Missing, "hermanliran@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D11"), "mailto:shmulmaor2@gmail.com", [Type]::Missing, [Type]::Missing, "shmulmaor2@gmail.com")

# Hyperlinks.Add() stamps its own "Hyperlink" look on the cell (new font/style) -
# put the plain e-mail-column formatting back, same as the other hyperlinked cells.
$ws.Range("C9").Copy()
$ws.Range("C10").PasteSpecial(-4122)
$ws.Range("D9").Copy()
$ws.Range("D11").PasteSpecial(-4122)

# --- Restore the view: scrolled back to A1, selection on the newly added rows ---
$ws.Activate()
$ws.Range("C11:D12").Select()
